{"js": "// Remove the \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" footer line, and\n// the blank paragraph that precedes them (the three paragraphs that sit\n// right after \"LOT2008: Bioqu\u00edmica II (Indica\u00e7\u00e3o de Conjunto)\" and right\n// before the trailing blank paragraph / page-break paragraph at the end\n// of the document).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOT2008: ...\") so the deletion is robust\n// even if surrounding content shifts.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOT2008\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Anchor paragraph 'LOT2008: ...' not found\");\n}\n\n// The three paragraphs to delete immediately follow the anchor:\n//   anchor + 1 -> empty paragraph\n//   anchor + 2 -> \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   anchor + 3 -> \"\u00a9 2020 . Contact: ...\"\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 3 && i < items.length; i++) {\n  toDelete.push(items[i]);\n}\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" footer line, and\n# the blank paragraph that precedes them (the three paragraphs that sit\n# right after \"LOT2008: Bioqu\u00edmica II (Indica\u00e7\u00e3o de Conjunto)\" and right\n# before the trailing blank paragraph / page-break paragraph at the end\n# of the document).\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# Locate the anchor paragraph (\"LOT2008: ...\") so the deletion is robust\n# even if surrounding content shifts.\n$anchorIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    if ($paras.Item($i).Range.Text -like \"*LOT2008*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph 'LOT2008: ...' not found\"\n}\n\n# The three paragraphs to delete immediately follow the anchor:\n#   anchor + 1 -> empty paragraph\n#   anchor + 2 -> \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   anchor + 3 -> \"\u00a9 2020 . Contact: ...\"\n# Delete from the last one back to the first so indices stay valid.\nfor ($offset = 3; $offset -ge 1; $offset--) {\n    $target = $d.Paragraphs.Item($anchorIndex + $offset)\n    $target.Range.Delete()\n}\n"}
